$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels in B1/C1
$ws.Range("B1").Value = "Retention Time"
$ws.Range("C1").Value = "Area"

# Update data rows: Vial id, Area value, Retention Time value
$ws.Range("A2").Value = "1:A,1"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 100

$ws.Range("A3").Value = "1:A,1"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 200

$ws.Range("A4").Value = "1:B,1"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 300

$ws.Range("A5").Value = "1:B,1"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 400
